# Adds the new "TGmonitor7" worksheet (sheet38) at the end of the workbook,
# fills it with the monitor spec content, updates the active-sheet/selection
# bookkeeping on the two sheets whose selection moved, and makes the new
# sheet the active tab - matching the target diff.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Selection housekeeping on the sheets that lose the "active" status.
#    (Do this before creating/activating the new sheet so it ends up last.)
# ---------------------------------------------------------------------------
$wsTGZS40014_30 = $wb.Worksheets.Item("TGZ-S-400-14_30")
$wsTGZS40014_30.Range("C19").Select()

$wsTGMmini = $wb.Worksheets.Item("TGMmini")
$wsTGMmini.Range("A1").Select()

# ---------------------------------------------------------------------------
# 2) Create the new sheet after the last existing sheet (commonHW_AI).
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "TGmonitor7"

# Column widths (characters), matching sibling spec sheets.
$ws.Columns.Item(1).ColumnWidth = 42.22
$ws.Columns.Item(2).ColumnWidth = 45.61

# ---------------------------------------------------------------------------
# 3) Cell content.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "POWER SUPPLY"
$ws.Range("B1").Value = " "

$ws.Range("A2").Value = "Supply voltage "
$ws.Range("A2").WrapText = $true
$ws.Range("B2").Value = "24 V DC (± 20 %)"

$ws.Range("A3").Value = "Recommended PSU current "
$ws.Range("A3").WrapText = $true
$ws.Range("B3").Value = "min. 300 mA"

$ws.Range("A5").Value = "CONNECTORS"
$ws.Range("B5").Value = " "

$ws.Range("A6").Value = "USB"
$ws.Range("B6").Value = "4 x USB 2.0, microUSB"

$ws.Range("A7").Value = "HDMI"
$ws.Range("B7").Value = "standard A"

$ws.Range("A8").Value = "Power"
$ws.Range("B8").Value = "1 x 4pin WEIDMÜLLER BLF 2.50/04/180 SN BK BX"

# Row 9 is a blank spacer row, but still materialised with empty styled cells.
$ws.Range("A9").NumberFormat = "General"
$ws.Range("B9").NumberFormat = "General"

$ws.Range("A10").Value = "DISPLAY"
$ws.Range("B10").Value = " "

$ws.Range("A11").Value = "Size"
$ws.Range("A11").WrapText = $true
$ws.Range("B11").Value = "7 inches"
$ws.Range("B11").WrapText = $true

$ws.Range("A12").Value = "Viewing Angles"
$ws.Range("A12").WrapText = $true
$ws.Range("B12").Value = "170°"
$ws.Range("B12").WrapText = $true

$ws.Range("A13").Value = "Resolution"
$ws.Range("A13").WrapText = $true
$ws.Range("B13").Value = "1024 x 600 pixels"
$ws.Range("B13").WrapText = $true
$ws.Rows.Item(13).RowHeight = 12.8

$ws.Range("A14").Value = "Display Area Dimensions"
$ws.Range("A14").WrapText = $true
$ws.Range("B14").Value = "154.21 x 85.92 mm"
$ws.Range("B14").WrapText = $true

$ws.Range("A15").Value = "Pixel Pitch"
$ws.Range("A15").WrapText = $true
$ws.Range("B15").Value = "150.6 µm (H) x 143.2 µm (V)"
$ws.Range("B15").WrapText = $true

$ws.Range("A16").Value = "Color Gamut"
$ws.Range("A16").WrapText = $true
$ws.Range("B16").Value = "45% NTSC"
$ws.Range("B16").WrapText = $true
$ws.Range("B16").NumberFormat = "0.00\ %"

$ws.Range("A17").Value = "Maximum Brightness"
$ws.Range("A17").WrapText = $true
$ws.Range("B17").Value = "300 cd/m²"
$ws.Range("B17").WrapText = $true
$ws.Rows.Item(17).RowHeight = 12.8

$ws.Range("A18").Value = "Contrast Ratio"
$ws.Range("A18").WrapText = $true
$ws.Range("B18").Value = "800:1"
$ws.Range("B18").WrapText = $true
$ws.Range("B18").NumberFormat = "@"
$ws.Rows.Item(18).RowHeight = 12.8

# ---------------------------------------------------------------------------
# 4) Page setup / print options / header-footer, matching sibling sheets.
# ---------------------------------------------------------------------------
$ws.PageSetup.LeftMargin = 56.7
$ws.PageSetup.RightMargin = 56.7
$ws.PageSetup.TopMargin = 75.8
$ws.PageSetup.BottomMargin = 75.8
$ws.PageSetup.HeaderMargin = 56.7
$ws.PageSetup.FooterMargin = 56.7

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Zoom = 100
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1
$ws.PageSetup.Orientation = 1

$ws.PageSetup.CenterHeader = "&""Times New Roman,obyčejné""&12&A"
$ws.PageSetup.CenterFooter = "&""Times New Roman,obyčejné""&12Stránka &P"

# ---------------------------------------------------------------------------
# 5) Activate the new sheet / selection last so it becomes the active tab.
# ---------------------------------------------------------------------------
$ws.Range("I20").Select()
